# Update "want to go" counts (column F) on several sheets to reflect the
# latest generated data (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 298
$ws1.Range("F5").Value = 349
$ws1.Range("F6").Value = 314
$ws1.Range("F9").Value = 742
$ws1.Range("F10").Value = 1561
$ws1.Range("F14").Value = 145
$ws1.Range("F16").Value = 2059

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F10").Value = 87

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F12").Value = 298
$ws4.Range("F14").Value = 349
$ws4.Range("F18").Value = 314
$ws4.Range("F22").Value = 87
$ws4.Range("F25").Value = 742
$ws4.Range("F26").Value = 1561
$ws4.Range("F31").Value = 145
$ws4.Range("F34").Value = 2059
